# Add new sample rows (76-86) to the "Tabelle1" worksheet, matching the
# new shared-string entries and the existing table's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$names = @(
  "HUC_squeez",
  "HUC_squeez2",
  "EGBA_kickles",
  "S3PM_kickakzent",
  "S3PM_kickakzent2",
  "EGBA_softkick",
  "DD_kiks0",
  "DD_kiks1",
  "DD_kiks2",
  "DD_kiks3",
  "JDP_kicky"
)

$firstNewRow = 76
$lastNewRow = $firstNewRow + $names.Count - 1

# Re-use the formatting (thin border around each cell) that's already used
# by the rest of the sample table, by copying it down from the last
# existing row into the newly added rows.
$ws.Range("A75:B75").Copy()
$ws.Range("A" + $firstNewRow + ":B" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $firstNewRow + $i
  $ws.Cells.Item($row, 1).Value = 75 + $i
  $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Leave the view roughly where the author left it.
$ws.Range("E64").Select()
